# Scheduled market-price refresh: updates currentAveragePrice(NQ/HQ) and
# the derived LevePrice/LeveProfit columns (H..N) across the per-job
# leve-profit tables (one per worksheet/job abbreviation).
# Column layout (1-based): H=8 currentAveragePrice, I=9 currentAveragePriceNQ,
# J=10 currentAveragePriceHQ, K=11 LevePriceNQ, L=12 LevePriceHQ,
# M=13 LeveProfitNQ, N=14 LeveProfitHQ.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(48, 8).Value = 1666.3334
$ws.Cells.Item(48, 9).Value = 200
$ws.Cells.Item(48, 10).Value = 2399.5
$ws.Cells.Item(48, 11).Value = 600
$ws.Cells.Item(48, 12).Value = 7198.5
$ws.Cells.Item(48, 13).Value = -308
$ws.Cells.Item(48, 14).Value = -7782.5

$ws.Cells.Item(56, 8).Value = 1666.3334
$ws.Cells.Item(56, 9).Value = 200
$ws.Cells.Item(56, 10).Value = 2399.5
$ws.Cells.Item(56, 11).Value = 600
$ws.Cells.Item(56, 12).Value = 7198.5
$ws.Cells.Item(56, 13).Value = -66
$ws.Cells.Item(56, 14).Value = -8266.5

$ws.Cells.Item(98, 8).Value = 9293.450000000001
$ws.Cells.Item(98, 9).Value = 9687.166999999999
$ws.Cells.Item(98, 11).Value = 9687.166999999999
$ws.Cells.Item(98, 13).Value = -8189.166999999999

$ws.Cells.Item(122, 8).Value = 9293.450000000001
$ws.Cells.Item(122, 9).Value = 9687.166999999999
$ws.Cells.Item(122, 11).Value = 29061.501
$ws.Cells.Item(122, 13).Value = -26611.501

$ws.Cells.Item(128, 8).Value = 37449.75
$ws.Cells.Item(128, 10).Value = 37449.75
$ws.Cells.Item(128, 12).Value = 37449.75
$ws.Cells.Item(128, 14).Value = -47409.75

$ws.Cells.Item(130, 8).Value = 35441
$ws.Cells.Item(130, 10).Value = 35441
$ws.Cells.Item(130, 12).Value = 35441
$ws.Cells.Item(130, 14).Value = -45481

$ws.Cells.Item(137, 8).Value = 960.2143
$ws.Cells.Item(137, 9).Value = 833.63635
$ws.Cells.Item(137, 10).Value = 1424.3334
$ws.Cells.Item(137, 11).Value = 2500.90905
$ws.Cells.Item(137, 12).Value = 4273.0002
$ws.Cells.Item(137, 13).Value = 49.09094999999979
$ws.Cells.Item(137, 14).Value = -9373.0002

$ws.Cells.Item(138, 8).Value = 2760.9119
$ws.Cells.Item(138, 9).Value = 1437.3182
$ws.Cells.Item(138, 10).Value = 5187.5
$ws.Cells.Item(138, 11).Value = 4311.9546
$ws.Cells.Item(138, 12).Value = 15562.5
$ws.Cells.Item(138, 13).Value = 828.0454
$ws.Cells.Item(138, 14).Value = -25842.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 399882.9
$ws.Cells.Item(32, 9).Value = 2586.6865
$ws.Cells.Item(32, 10).Value = 5723652.5
$ws.Cells.Item(32, 11).Value = 2586.6865
$ws.Cells.Item(32, 12).Value = 5723652.5
$ws.Cells.Item(32, 13).Value = -2299.6865
$ws.Cells.Item(32, 14).Value = -5724226.5

$ws.Cells.Item(74, 8).Value = 822.5
$ws.Cells.Item(74, 9).Value = 701.4211
$ws.Cells.Item(74, 10).Value = 1282.6
$ws.Cells.Item(74, 11).Value = 701.4211
$ws.Cells.Item(74, 12).Value = 1282.6
$ws.Cells.Item(74, 13).Value = 172.5789
$ws.Cells.Item(74, 14).Value = -3030.6

$ws.Cells.Item(77, 8).Value = 822.5
$ws.Cells.Item(77, 9).Value = 701.4211
$ws.Cells.Item(77, 10).Value = 1282.6
$ws.Cells.Item(77, 11).Value = 3507.1055
$ws.Cells.Item(77, 12).Value = 6413
$ws.Cells.Item(77, 13).Value = 860.8944999999999
$ws.Cells.Item(77, 14).Value = -15149

$ws.Cells.Item(132, 8).Value = 1025.1023
$ws.Cells.Item(132, 9).Value = 650.77466
$ws.Cells.Item(132, 10).Value = 2588.4707
$ws.Cells.Item(132, 11).Value = 1952.32398
$ws.Cells.Item(132, 12).Value = 7765.4121
$ws.Cells.Item(132, 13).Value = 577.6760199999999
$ws.Cells.Item(132, 14).Value = -12825.4121

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2400.516
$ws.Cells.Item(31, 9).Value = 2066
$ws.Cells.Item(31, 10).Value = 4140
$ws.Cells.Item(31, 11).Value = 2066
$ws.Cells.Item(31, 12).Value = 4140
$ws.Cells.Item(31, 13).Value = -1771
$ws.Cells.Item(31, 14).Value = -4730

$ws.Cells.Item(34, 8).Value = 2400.516
$ws.Cells.Item(34, 9).Value = 2066
$ws.Cells.Item(34, 10).Value = 4140
$ws.Cells.Item(34, 11).Value = 2066
$ws.Cells.Item(34, 12).Value = 4140
$ws.Cells.Item(34, 13).Value = -1864
$ws.Cells.Item(34, 14).Value = -4544

$ws.Cells.Item(132, 8).Value = 32995.78
$ws.Cells.Item(132, 9).Value = 878.087
$ws.Cells.Item(132, 10).Value = 115074.336
$ws.Cells.Item(132, 11).Value = 2634.261
$ws.Cells.Item(132, 12).Value = 345223.008
$ws.Cells.Item(132, 13).Value = -104.261
$ws.Cells.Item(132, 14).Value = -350283.008

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1081
$ws.Cells.Item(113, 9).Value = 994.4
$ws.Cells.Item(113, 10).Value = 1094.9678
$ws.Cells.Item(113, 11).Value = 2983.2
$ws.Cells.Item(113, 12).Value = 3284.9034
$ws.Cells.Item(113, 13).Value = -813.1999999999998
$ws.Cells.Item(113, 14).Value = -7624.903399999999

$ws.Cells.Item(122, 8).Value = 646.5454999999999
$ws.Cells.Item(122, 10).Value = 1076
$ws.Cells.Item(122, 12).Value = 9684
$ws.Cells.Item(122, 14).Value = -14584

$ws.Cells.Item(132, 8).Value = 2231.913
$ws.Cells.Item(132, 9).Value = 1292
$ws.Cells.Item(132, 10).Value = 2836.1428
$ws.Cells.Item(132, 11).Value = 11628
$ws.Cells.Item(132, 12).Value = 25525.2852
$ws.Cells.Item(132, 13).Value = -9098
$ws.Cells.Item(132, 14).Value = -30585.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 15035010
$ws.Cells.Item(11, 9).Value = 30060020
$ws.Cells.Item(11, 10).Value = 10000
$ws.Cells.Item(11, 11).Value = 30060020
$ws.Cells.Item(11, 12).Value = 10000
$ws.Cells.Item(11, 13).Value = -30059881
$ws.Cells.Item(11, 14).Value = -10278

$ws.Cells.Item(18, 8).Value = 7006
$ws.Cells.Item(18, 10).Value = 7006
$ws.Cells.Item(18, 12).Value = 7006
$ws.Cells.Item(18, 14).Value = -7592

$ws.Cells.Item(21, 8).Value = 456409.1
$ws.Cells.Item(21, 9).Value = 5000500
$ws.Cells.Item(21, 11).Value = 5000500
$ws.Cells.Item(21, 13).Value = -5000327

$ws.Cells.Item(30, 8).Value = 456409.1
$ws.Cells.Item(30, 9).Value = 5000500
$ws.Cells.Item(30, 11).Value = 5000500
$ws.Cells.Item(30, 13).Value = -5000395

$ws.Cells.Item(33, 8).Value = 111556744
$ws.Cells.Item(33, 10).Value = 111556744
$ws.Cells.Item(33, 12).Value = 111556744
$ws.Cells.Item(33, 14).Value = -111557248

$ws.Cells.Item(35, 8).Value = 5000
$ws.Cells.Item(35, 10).Value = 5000
$ws.Cells.Item(35, 12).Value = 5000
$ws.Cells.Item(35, 14).Value = -5596

$ws.Cells.Item(43, 8).Value = 39514.145
$ws.Cells.Item(43, 9).Value = 999.5
$ws.Cells.Item(43, 10).Value = 54920
$ws.Cells.Item(43, 11).Value = 999.5
$ws.Cells.Item(43, 12).Value = 54920
$ws.Cells.Item(43, 13).Value = -848.5
$ws.Cells.Item(43, 14).Value = -55222

$ws.Cells.Item(107, 8).Value = 287.25
$ws.Cells.Item(107, 9).Value = 204.3077
$ws.Cells.Item(107, 10).Value = 441.2857
$ws.Cells.Item(107, 11).Value = 204.3077
$ws.Cells.Item(107, 12).Value = 441.2857
$ws.Cells.Item(107, 13).Value = 1715.6923
$ws.Cells.Item(107, 14).Value = -4281.2857

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 2500
$ws.Cells.Item(20, 10).Value = 2500
$ws.Cells.Item(20, 12).Value = 2500
$ws.Cells.Item(20, 14).Value = -2952

$ws.Cells.Item(57, 8).Value = 49250
$ws.Cells.Item(57, 10).Value = 49250
$ws.Cells.Item(57, 12).Value = 49250
$ws.Cells.Item(57, 14).Value = -50382

$ws.Cells.Item(132, 8).Value = 3335.9863
$ws.Cells.Item(132, 9).Value = 5057.769
$ws.Cells.Item(132, 10).Value = 1361
$ws.Cells.Item(132, 11).Value = 15173.307
$ws.Cells.Item(132, 12).Value = 4083
$ws.Cells.Item(132, 13).Value = -12643.307
$ws.Cells.Item(132, 14).Value = -9143

$ws.Cells.Item(136, 8).Value = 4100.04
$ws.Cells.Item(136, 9).Value = 2644.2222
$ws.Cells.Item(136, 10).Value = 7843.5713
$ws.Cells.Item(136, 11).Value = 7932.6666
$ws.Cells.Item(136, 12).Value = 23530.7139
$ws.Cells.Item(136, 13).Value = -5382.6666
$ws.Cells.Item(136, 14).Value = -28630.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(37, 8).Value = 3599
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 11).Value = 0
# LeveProfitNQ (M37) no longer populated for this row
$ws.Cells.Item(37, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 19139416
$ws.Cells.Item(132, 9).Value = 25511262
$ws.Cells.Item(132, 10).Value = 1793838.4
$ws.Cells.Item(132, 11).Value = 76533786
$ws.Cells.Item(132, 12).Value = 5381515.199999999
$ws.Cells.Item(132, 13).Value = -76531256
$ws.Cells.Item(132, 14).Value = -5386575.199999999